# chore: update Sheets via scheduled runner
# Refreshes market-price-derived columns (currentAveragePrice*, LevePrice*,
# LeveProfit*) on several leve rows across the ALC/ARM/BSM/CRP/CUL/LTW
# sheets, matching a scheduled Golem_Profits price-data pull.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H3").Value2  = 60714.285
$ws.Range("J3").Value2  = 60714.285
$ws.Range("L3").Value2  = 60714.285
$ws.Range("N3").Value2  = -60942.285

$ws.Range("H31").Value2 = 102.5
$ws.Range("I31").Value2 = 102.5
$ws.Range("K31").Value2 = 307.5
$ws.Range("M31").Value2 = -77.5

$ws.Range("H33").Value2 = 61.136364
$ws.Range("I33").Value2 = 61.136364
$ws.Range("K33").Value2 = 61.136364
$ws.Range("M33").Value2 = 167.863636

$ws.Range("H75").Value2 = 45000
$ws.Range("J75").Value2 = 45000
$ws.Range("L75").Value2 = 45000
$ws.Range("N75").Value2 = -46872

$ws.Range("H78").Value2 = 45000
$ws.Range("J78").Value2 = 45000
$ws.Range("L78").Value2 = 135000
$ws.Range("N78").Value2 = -144360

$ws.Range("H86").Value2 = 10100.6
$ws.Range("I86").Value2 = 0
$ws.Range("J86").Value2 = 10100.6
$ws.Range("K86").Value2 = 0
$ws.Range("L86").Value2 = 10100.6
$ws.Range("M86").Value2 = $null
$ws.Range("N86").Value2 = -12346.6

$ws.Range("H89").Value2 = 10100.6
$ws.Range("I89").Value2 = 0
$ws.Range("J89").Value2 = 10100.6
$ws.Range("K89").Value2 = 0
$ws.Range("L89").Value2 = 50503
$ws.Range("M89").Value2 = $null
$ws.Range("N89").Value2 = -61735

$ws.Range("H96").Value2 = 725.5
$ws.Range("I96").Value2 = 356.5
$ws.Range("J96").Value2 = 971.5
$ws.Range("K96").Value2 = 1069.5
$ws.Range("L96").Value2 = 2914.5
$ws.Range("M96").Value2 = 303.5
$ws.Range("N96").Value2 = -5660.5

$ws.Range("H102").Value2 = 60714.285
$ws.Range("J102").Value2 = 60714.285
$ws.Range("L102").Value2 = 60714.285
$ws.Range("N102").Value2 = -67204.285

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H10").Value2 = 1426.25
$ws.Range("I10").Value2 = 1566.6666
$ws.Range("J10").Value2 = 1005
$ws.Range("K10").Value2 = 1566.6666
$ws.Range("L10").Value2 = 1005
$ws.Range("M10").Value2 = -1396.6666
$ws.Range("N10").Value2 = -1345

$ws.Range("H12").Value2 = 267.5
$ws.Range("I12").Value2 = 300
$ws.Range("J12").Value2 = 235
$ws.Range("K12").Value2 = 300
$ws.Range("L12").Value2 = 235
$ws.Range("M12").Value2 = -127
$ws.Range("N12").Value2 = -581

$ws.Range("H13").Value2 = 424.5
$ws.Range("J13").Value2 = 424.5
$ws.Range("L13").Value2 = 424.5
$ws.Range("N13").Value2 = -712.5

$ws.Range("H16").Value2 = 597
$ws.Range("I16").Value2 = 1000
$ws.Range("J16").Value2 = 194
$ws.Range("K16").Value2 = 1000
$ws.Range("L16").Value2 = 194
$ws.Range("M16").Value2 = -713
$ws.Range("N16").Value2 = -768

$ws.Range("H30").Value2 = 850
$ws.Range("I30").Value2 = 850
$ws.Range("K30").Value2 = 850
$ws.Range("M30").Value2 = -700

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H26").Value2 = 26468.5
$ws.Range("I26").Value2 = 26468.5
$ws.Range("K26").Value2 = 26468.5
$ws.Range("M26").Value2 = -26176.5

$ws.Range("H105").Value2 = 1424.8334
$ws.Range("J105").Value2 = 2400
$ws.Range("L105").Value2 = 2400
$ws.Range("N105").Value2 = -5894

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H86").Value2 = 8182
$ws.Range("I86").Value2 = 8576.875
$ws.Range("J86").Value2 = 7129
$ws.Range("K86").Value2 = 8576.875
$ws.Range("L86").Value2 = 7129
$ws.Range("M86").Value2 = -7453.875
$ws.Range("N86").Value2 = -9375

$ws.Range("H89").Value2 = 8182
$ws.Range("I89").Value2 = 8576.875
$ws.Range("J89").Value2 = 7129
$ws.Range("K89").Value2 = 42884.375
$ws.Range("L89").Value2 = 35645
$ws.Range("M89").Value2 = -37268.375
$ws.Range("N89").Value2 = -46877

$ws.Range("H107").Value2 = 226.88889
$ws.Range("I107").Value2 = 242.8125
$ws.Range("K107").Value2 = 242.8125
$ws.Range("M107").Value2 = 1677.1875

$ws.Range("H132").Value2 = 475
$ws.Range("I132").Value2 = 450
$ws.Range("J132").Value2 = 550
$ws.Range("K132").Value2 = 1350
$ws.Range("L132").Value2 = 1650
$ws.Range("M132").Value2 = 1180
$ws.Range("N132").Value2 = -6710

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H36").Value2 = 20
$ws.Range("I36").Value2 = 20
$ws.Range("K36").Value2 = 60
$ws.Range("M36").Value2 = 109

$ws.Range("H64").Value2 = 12
$ws.Range("I64").Value2 = 12
$ws.Range("K64").Value2 = 36
$ws.Range("M64").Value2 = 234

$ws.Range("H67").Value2 = 12
$ws.Range("I67").Value2 = 12
$ws.Range("K67").Value2 = 36
$ws.Range("M67").Value2 = 900

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H7").Value2 = 0
$ws.Range("I7").Value2 = 0
$ws.Range("K7").Value2 = 0
$ws.Range("M7").Value2 = $null

$ws.Range("H22").Value2 = 2279.2307
$ws.Range("I22").Value2 = 999.8570999999999
$ws.Range("J22").Value2 = 3771.8333
$ws.Range("K22").Value2 = 999.8570999999999
$ws.Range("L22").Value2 = 3771.8333
$ws.Range("M22").Value2 = -704.8570999999999
$ws.Range("N22").Value2 = -4361.8333

$ws.Range("H27").Value2 = 2279.2307
$ws.Range("I27").Value2 = 999.8570999999999
$ws.Range("J27").Value2 = 3771.8333
$ws.Range("K27").Value2 = 999.8570999999999
$ws.Range("L27").Value2 = 3771.8333
$ws.Range("M27").Value2 = -892.8570999999999
$ws.Range("N27").Value2 = -3985.8333

$ws.Range("H31").Value2 = 4416.2856
$ws.Range("I31").Value2 = 1507.5
$ws.Range("J31").Value2 = 5579.8
$ws.Range("K31").Value2 = 1507.5
$ws.Range("L31").Value2 = 5579.8
$ws.Range("M31").Value2 = -1259.5
$ws.Range("N31").Value2 = -6075.8

$ws.Range("H40").Value2 = 1730001.4
$ws.Range("I40").Value2 = 95002
$ws.Range("J40").Value2 = 5000000
$ws.Range("K40").Value2 = 95002
$ws.Range("L40").Value2 = 5000000
$ws.Range("M40").Value2 = -94866
$ws.Range("N40").Value2 = -5000272

$ws.Range("H61").Value2 = 9666.666999999999
$ws.Range("I61").Value2 = 9666.666999999999
$ws.Range("K61").Value2 = 9666.666999999999
$ws.Range("M61").Value2 = -9464.666999999999

$ws.Range("H113").Value2 = 9666.666999999999
$ws.Range("I113").Value2 = 9666.666999999999
$ws.Range("K113").Value2 = 9666.666999999999
$ws.Range("M113").Value2 = -7496.666999999999

$ws.Range("H122").Value2 = 3496.25
$ws.Range("J122").Value2 = 3496.25
$ws.Range("L122").Value2 = 10488.75
$ws.Range("N122").Value2 = -15388.75

$ws.Range("H126").Value2 = 0
$ws.Range("I126").Value2 = 0
$ws.Range("K126").Value2 = 0
$ws.Range("M126").Value2 = $null
